{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Summary of the target edit:\n//  1. The document title paragraph (\"Builder Pattern Example with Computer\n//     Database\") is split into two bold, 16pt paragraphs: a new\n//     \"Exercise-3 \" paragraph followed by the original title text.\n//  2. A handful of paragraphs throughout the body have stray markdown\n//     backticks (`) removed from around inline code terms (Computer,\n//     ComputerDatabase, main, toString(), build()), and two list items\n//     lose their leading \"-\" bullet character/prefix.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// --- Step 1: plain text clean-ups (done first, while paragraph indices\n// still match the original document layout -- the title split below\n// would otherwise shift every later index by one). Each entry is matched\n// by its exact current text so the edit is resilient to being re-run in\n// a slightly different order.\nconst textFixes = [\n  [\n    \"The `Computer` class represents a computer with various attributes:\",\n    \"The Computer class represents a computer with various attributes:\",\n  ],\n  [\n    \"- Getter methods for all attributes\",\n    \" Getter methods for all attributes\",\n  ],\n  [\n    \"- A `toString()` method for string representation\",\n    \"toString() method for string representation\",\n  ],\n  [\n    \"The `Computer` class contains an inner `Builder` class that:\",\n    \"The Computer class contains an inner `Builder` class that:\",\n  ],\n  [\n    \"- Includes a `build()` method to create the final Computer object\",\n    \"- Includes a build() method to create the final Computer object\",\n  ],\n  [\n    \"The `ComputerDatabase` class serves as a simple in-memory database for storing Computer objects.\",\n    \"The ComputerDatabase class serves as a simple in-memory database for storing Computer objects.\",\n  ],\n  [\n    \"This class contains the `main` method and serves as the entry point of the application. It demonstrates:\",\n    \"This class contains the main method and serves as the entry point of the application. It demonstrates:\",\n  ],\n  [\n    \"   - Implement the `build()` method to create and return a new Computer object\",\n    \"   - Implement the build() method to create and return a new Computer object\",\n  ],\n];\n\nfor (const [oldText, newText] of textFixes) {\n  const match = paragraphs.items.find((p) => p.text === oldText);\n  if (match) {\n    match.insertText(newText, \"Replace\");\n  }\n}\nawait context.sync();\n\n// --- Step 2: split the title paragraph into \"Exercise-3 \" + the\n// original title, both bold and 16pt (sz/szCs 32 half-points).\nconst titleParagraph = body.paragraphs.items[0];\ntitleParagraph.load(\"text\");\nawait context.sync();\n\ntitleParagraph.insertParagraph(\"Exercise-3 \", \"Before\");\nawait context.sync();\n\n// Re-fetch the paragraphs after the insertion -- applying formatting to\n// the proxy objects returned before the sync above does not reliably\n// stick, so grab fresh Paragraph objects by (now-shifted) index.\nconst refreshedParagraphs = context.document.body.paragraphs;\nrefreshedParagraphs.load(\"text\");\nawait context.sync();\n\nconst exercisePara = refreshedParagraphs.items[0];\nconst titlePara = refreshedParagraphs.items[1];\n\nfor (const para of [exercisePara, titlePara]) {\n  para.font.bold = true;\n  para.font.boldBidirectional = true;\n  para.font.size = 16;\n  para.font.sizeBidirectional = 16;\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Summary of the target edit:\n#  1. The document title paragraph (\"Builder Pattern Example with Computer\n#     Database\") is split into two bold, 16pt paragraphs: a new\n#     \"Exercise-3 \" paragraph followed by the original title text.\n#  2. A handful of paragraphs throughout the body have stray markdown\n#     backticks (`) removed from around inline code terms (Computer,\n#     ComputerDatabase, main, toString(), build()), and two list items\n#     lose their leading \"-\" bullet character/prefix.\n\n$d = $word.ActiveDocument\n\nfunction Replace-ParagraphText($oldText, $newText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n\n# --- Step 1: plain text clean-ups (done first, while paragraph indices\n# still match the original document layout -- the title split below\n# would otherwise shift every later paragraph index by one). Each\n# original string is unique in the document, so a single Find/Replace\n# pass touches exactly the right paragraph.\nReplace-ParagraphText 'The `Computer` class represents a computer with various attributes:' 'The Computer class represents a computer with various attributes:'\nReplace-ParagraphText '- Getter methods for all attributes' ' Getter methods for all attributes'\nReplace-ParagraphText '- A `toString()` method for string representation' 'toString() method for string representation'\nReplace-ParagraphText 'The `Computer` class contains an inner `Builder` class that:' 'The Computer class contains an inner `Builder` class that:'\nReplace-ParagraphText '- Includes a `build()` method to create the final Computer object' '- Includes a build() method to create the final Computer object'\nReplace-ParagraphText 'The `ComputerDatabase` class serves as a simple in-memory database for storing Computer objects.' 'The ComputerDatabase class serves as a simple in-memory database for storing Computer objects.'\nReplace-ParagraphText 'This class contains the `main` method and serves as the entry point of the application. It demonstrates:' 'This class contains the main method and serves as the entry point of the application. It demonstrates:'\nReplace-ParagraphText '   - Implement the `build()` method to create and return a new Computer object' '   - Implement the build() method to create and return a new Computer object'\n\n# --- Step 2: split the title paragraph into \"Exercise-3 \" + the\n# original title, both bold and 16pt (sz/szCs 32 half-points).\n$titlePara = $d.Paragraphs.Item(1)\n$titleRange = $titlePara.Range\n$titleRange.Collapse(1)\n$titleRange.InsertBefore(\"Exercise-3 `r\")\n\n$exercisePara = $d.Paragraphs.Item(1)\n$exercisePara.Range.Font.Bold = 1\n$exercisePara.Range.Font.BoldBi = 1\n$exercisePara.Range.Font.Size = 16\n$exercisePara.Range.Font.SizeBi = 16\n\n$titlePara2 = $d.Paragraphs.Item(2)\n$titlePara2.Range.Font.Bold = 1\n$titlePara2.Range.Font.BoldBi = 1\n$titlePara2.Range.Font.Size = 16\n$titlePara2.Range.Font.SizeBi = 16\n"}
